$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("D1").Value = "regu"
$ws.Range("E1").Value = "kelompok_id"
$ws.Range("F1").Value = "desa"

# Update row 2 (Udin)
$ws.Range("B2").Value = 111
$ws.Range("D2").Value = "Maroon"
$ws.Range("E2").Value = "KM7"
$ws.Range("F2").Value = "Batam"

# Update row 3 (Maimunah)
$ws.Range("B3").Value = 112
$ws.Range("D3").Value = "Biru"
$ws.Range("E3").Value = "km7"
$ws.Range("F3").Value = "batam"

# Update selection to F3
$ws.Range("F3").Select()
